# Regenerate merged AHB files
# - Rename the "_old" / "_new" suffixed header strings to "_FV2210" / "_FV2304"
# - Wrap the data range in a native Excel Table ("Table1")
# - Freeze the header row (split/freeze pane below row 1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename shared header strings (old -> FV2210, new -> FV2304)
$ws.Cells.Replace("_old", "_FV2210")
$ws.Cells.Replace("_new", "_FV2304")

# 2) Turn the used range into an Excel table so the header row becomes
#    structured columns (names are picked up from row 1, which we just renamed).
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U62"), $null, 1)
$tbl.Name = "Table1"

# 3) Freeze the top (header) row.
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
